# Update cryptos list data (price + 1h volume change) per latest scrape.
# D-column price cells are stored as plain text in the source data (e.g.
# "64.336.37", "0.0000187"), so we force text format before assigning the
# new values to avoid Excel auto-converting them to numbers/scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.481.65"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.418.34"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.76"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.17"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  +5.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.423.45"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.008.86"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  -3.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.92"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.537.56"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.411.41"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.02"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.11"
$ws.Range("E21").Value = "  -4.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.97"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.557"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.41"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("E27").Value = "  +8.23%  "
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  +3.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.02"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.11"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.22"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("E35").Value = "  +6.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.98"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0761"
$ws.Range("E38").Value = "  -1.83%  "

# Rows 39/40: RenderToken and EnergySwap swapped positions, with updated data
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.83"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.76"
$ws.Range("E40").Value = "  -3.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.846.82"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.87"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.44"
$ws.Range("E44").Value = "  +8.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0314"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.769"
$ws.Range("E46").Value = "  -0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.11"
$ws.Range("E47").Value = "  +6.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.08"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.59"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.858"
$ws.Range("E51").Value = "  -2.01%  "
